$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update example column (C) for existing concept rows, and move/replace rows
# so the table matches the new spec examples.

$ws.Range("A2").Value = "cell"
$ws.Range("B2").Value = "Single cell in the sheet"
$ws.Range("C2").Value = "let A1 = CELL(5);"

$ws.Range("A3").Value = "range"
$ws.Range("B3").Value = "Range of cells"
$ws.Range("C3").Value = "let rng = 1:5~1;"

$ws.Range("A4").Value = "array"
$ws.Range("B4").Value = "List of values of same primitive data type"
$ws.Range("C4").Value = "let num = [1,2,3,4,5];"

$ws.Range("A5").Value = "table"
$ws.Range("B5").Value = "Represents structured range of cells with headers and data rows"
$ws.Range("C5").Value = "let tab1 = TABLE();"

$ws.Range("A6").Value = "formula"
$ws.Range("B6").Value = "Datatype that holds formula"
$ws.Range("C6").Value = "let f1 = FORMULA([a, b], SUM);"

$ws.Range("A7").Value = "class"
$ws.Range("B7").Value = "User-defined datatype combining multiple primitive and/or non-primitive data types"
$ws.Range("C7").Value = "struct example {`n    let id = 1, speaker = `"john`";`n    let country = `"au`";`n};`nconst exampleObj = example(101, `"widget`", 19.99);"

# Column width adjustment for column C
# (Excel quantizes ColumnWidth to whole-pixel increments for the sheet's
# default font, so we pick the input that lands on the closest achievable
# stored width to the target 28.36328125.)
$ws.Columns.Item(3).ColumnWidth = 27.42

# Row height adjustments (row 4 keeps its default/natural height, so it is
# intentionally left untouched here)
$ws.Rows.Item(2).RowHeight = 21
$ws.Rows.Item(3).RowHeight = 18.5
$ws.Rows.Item(5).RowHeight = 21.5
$ws.Rows.Item(6).RowHeight = 20.5
$ws.Rows.Item(7).RowHeight = 84

# Update the active selection
$ws.Range("D10").Select()
